$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before A (shifts existing A:G data to B:H)
[void]$ws.Columns.Item(1).Insert()

# Insert a new row before 1 (shifts existing rows 1:9 data to 2:10)
[void]$ws.Rows.Item(1).Insert()

# New label column (column A) for the existing parameter rows
$ws.Cells.Item(3, 1).Value = "center"
$ws.Cells.Item(4, 1).Value = "sigma"
$ws.Cells.Item(5, 1).Value = "sigma min"
$ws.Cells.Item(6, 1).Value = "amplitude"
$ws.Cells.Item(7, 1).Value = "ampl. Min"

# New header row (row 1): names for each peak's parameter column
$ws.Cells.Item(1, 2).Value = "D4"
$ws.Cells.Item(1, 3).Value = "D5"
$ws.Cells.Item(1, 5).Value = "D3a"
$ws.Cells.Item(1, 6).Value = "D3b"
$ws.Cells.Item(1, 4).Value = "D1"
$ws.Cells.Item(1, 7).Value = "G"
$ws.Cells.Item(1, 8).Value = "D2"

$ws.Cells.Item(1, 1).Value = "name"
$ws.Cells.Item(2, 1).Value = "activate peak"

# New label rows for the appended fraction parameters
$ws.Cells.Item(8, 1).Value = "fraction"
$ws.Cells.Item(9, 1).Value = "fraction min"
$ws.Cells.Item(10, 1).Value = "fraction max"

# Widen column A so the labels are readable
$ws.Columns.Item(1).ColumnWidth = 11.33

# Select the header row, matching the authored selection state
[void]$ws.Range("A1:H1").Select()
